# Scheduled runner update: refresh cached Universalis market-price snapshots
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ /
# LevePriceHQ / LeveProfitNQ / LeveProfitHQ columns) across the ALC, ARM, BSM, CRP,
# CUL, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H40").Value = 1436.3636
$ws.Range("I40").Value = 1325
$ws.Range("K40").Value = 1325
$ws.Range("M40").Value = -1150

$ws.Range("H98").Value = 1061.5883
$ws.Range("I98").Value = 969.8
$ws.Range("J98").Value = 1750
$ws.Range("K98").Value = 969.8
$ws.Range("L98").Value = 1750
$ws.Range("M98").Value = 528.2
$ws.Range("N98").Value = -4746

$ws.Range("H111").Value = 5605.6
$ws.Range("I111").Value = 5363.5713
$ws.Range("J111").Value = 5817.375
$ws.Range("K111").Value = 16090.7139
$ws.Range("L111").Value = 17452.125
$ws.Range("M111").Value = -13023.7139
$ws.Range("N111").Value = -23586.125

$ws.Range("H112").Value = 7768
$ws.Range("I112").Value = 1100
$ws.Range("J112").Value = 8339.543
$ws.Range("K112").Value = 3300
$ws.Range("L112").Value = 25018.629
$ws.Range("M112").Value = -2192
$ws.Range("N112").Value = -27234.629

$ws.Range("H122").Value = 1061.5883
$ws.Range("I122").Value = 969.8
$ws.Range("J122").Value = 1750
$ws.Range("K122").Value = 2909.4
$ws.Range("L122").Value = 5250
$ws.Range("M122").Value = -459.3999999999996
$ws.Range("N122").Value = -10150

$ws.Range("H129").Value = 1212.3448
$ws.Range("I129").Value = 562.7273
$ws.Range("J129").Value = 1609.3334
$ws.Range("K129").Value = 1688.1819
$ws.Range("L129").Value = 4828.0002
$ws.Range("M129").Value = 3311.8181
$ws.Range("N129").Value = -14828.0002

$ws.Range("H132").Value = 4474.7666
$ws.Range("I132").Value = 4284.276
$ws.Range("K132").Value = 12852.828
$ws.Range("M132").Value = -10322.828

$ws.Range("H138").Value = 180245.67
$ws.Range("J138").Value = 280692.78
$ws.Range("L138").Value = 842078.3400000001
$ws.Range("N138").Value = -852358.3400000001

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H61").Value = 3038.4614
$ws.Range("I61").Value = 2266.6667
$ws.Range("J61").Value = 3700
$ws.Range("K61").Value = 2266.6667
$ws.Range("L61").Value = 3700
$ws.Range("M61").Value = -2054.6667
$ws.Range("N61").Value = -4124

$ws.Range("H68").Value = 56819.8
$ws.Range("J68").Value = 56819.8
$ws.Range("L68").Value = 56819.8
$ws.Range("N68").Value = -58441.8

$ws.Range("H71").Value = 56819.8
$ws.Range("J71").Value = 56819.8
$ws.Range("L71").Value = 170459.4
$ws.Range("N71").Value = -178571.4

$ws.Range("H74").Value = 1011.17145
$ws.Range("I74").Value = 985.5417
$ws.Range("K74").Value = 985.5417
$ws.Range("M74").Value = -111.5417

$ws.Range("H77").Value = 1011.17145
$ws.Range("I77").Value = 985.5417
$ws.Range("K77").Value = 4927.7085
$ws.Range("M77").Value = -559.7084999999997

$ws.Range("H122").Value = 2065.3794
$ws.Range("I122").Value = 1733.3914
$ws.Range("K122").Value = 5200.174199999999
$ws.Range("M122").Value = -2750.174199999999

$ws.Range("H136").Value = 3038.4614
$ws.Range("I136").Value = 2266.6667
$ws.Range("J136").Value = 3700
$ws.Range("K136").Value = 6800.000100000001
$ws.Range("L136").Value = 11100
$ws.Range("M136").Value = -4250.000100000001
$ws.Range("N136").Value = -16200

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H22").Value = 2545.4644
$ws.Range("I22").Value = 2545.4644
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 2545.4644
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -2372.4644
$ws.Range("N22").ClearContents()

$ws.Range("H94").Value = 1000
$ws.Range("I94").Value = 1000
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 1000
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = -549
$ws.Range("N94").Value = -1902

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H16").Value = 827.2
$ws.Range("I16").Value = 770.3333
$ws.Range("J16").Value = 912.5
$ws.Range("K16").Value = 770.3333
$ws.Range("L16").Value = 912.5
$ws.Range("M16").Value = -483.3333
$ws.Range("N16").Value = -1486.5

$ws.Range("H31").Value = 2301.9487
$ws.Range("I31").Value = 944.1539
$ws.Range("J31").Value = 5017.5386
$ws.Range("K31").Value = 944.1539
$ws.Range("L31").Value = 5017.5386
$ws.Range("M31").Value = -649.1539
$ws.Range("N31").Value = -5607.5386

$ws.Range("H34").Value = 2301.9487
$ws.Range("I34").Value = 944.1539
$ws.Range("J34").Value = 5017.5386
$ws.Range("K34").Value = 944.1539
$ws.Range("L34").Value = 5017.5386
$ws.Range("M34").Value = -742.1539
$ws.Range("N34").Value = -5421.5386

$ws.Range("H113").Value = 827.2
$ws.Range("I113").Value = 770.3333
$ws.Range("J113").Value = 912.5
$ws.Range("K113").Value = 770.3333
$ws.Range("L113").Value = 912.5
$ws.Range("M113").Value = 1399.6667
$ws.Range("N113").Value = -5252.5

$ws.Range("H134").Value = 1597.5333
$ws.Range("I134").Value = 1637.4166
$ws.Range("J134").Value = 1438
$ws.Range("K134").Value = 4912.2498
$ws.Range("L134").Value = 4314
$ws.Range("M134").Value = -2377.2498
$ws.Range("N134").Value = -9384

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H122").Value = 5448.143
$ws.Range("I122").Value = 494.33334
$ws.Range("J122").Value = 17832.666
$ws.Range("K122").Value = 4449.00006
$ws.Range("L122").Value = 160493.994
$ws.Range("M122").Value = -1999.00006
$ws.Range("N122").Value = -165393.994

$ws.Range("H131").Value = 955.5443
$ws.Range("J131").Value = 1009.9726
$ws.Range("L131").Value = 3029.9178
$ws.Range("N131").Value = -13109.9178

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H7").Value = 71431096
$ws.Range("I7").Value = 111112500
$ws.Range("J7").Value = 4562
$ws.Range("K7").Value = 111112500
$ws.Range("L7").Value = 4562
$ws.Range("M7").Value = -111112388
$ws.Range("N7").Value = -4786

$ws.Range("H40").Value = 38465940
$ws.Range("I40").Value = 76926200
$ws.Range("J40").Value = 5677.3076
$ws.Range("K40").Value = 76926200
$ws.Range("L40").Value = 5677.3076
$ws.Range("M40").Value = -76926064
$ws.Range("N40").Value = -5949.3076

$ws.Range("H116").Value = 35000
$ws.Range("J116").Value = 35000
$ws.Range("L116").Value = 35000
$ws.Range("N116").Value = -44178

$ws.Range("H126").Value = 71431096
$ws.Range("I126").Value = 111112500
$ws.Range("J126").Value = 4562
$ws.Range("K126").Value = 333337500
$ws.Range("L126").Value = 13686
$ws.Range("M126").Value = -333335030
$ws.Range("N126").Value = -18626

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H133").Value = 40715
$ws.Range("J133").Value = 40715
$ws.Range("L133").Value = 40715
$ws.Range("N133").Value = -50835

$ws.Range("H136").Value = 3394.56
$ws.Range("I136").Value = 2704.3333
$ws.Range("J136").Value = 4429.9
$ws.Range("K136").Value = 8112.999899999999
$ws.Range("L136").Value = 13289.7
$ws.Range("M136").Value = -5562.999899999999
$ws.Range("N136").Value = -18389.7
